$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.060.21"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.617.62"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.84%  "
$ws.Range("D9").Value = "2.620.33"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  -5.60%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "3.075.06"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "60.074.25"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").Value = "2.619.01"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.993"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("D28").Value = "0.0₃0805"
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("E33").Value = "  -4.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  -5.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.916"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -5.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.860"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("E40").Value = "  -4.14%  "
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0545"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("D51").Value = "1.959.00"
$ws.Range("E51").Value = "  +0.05%  "
